$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "sequences/278857_percepttesting_sequence_00.csv"
$ws.Range("A3").Value = "sequences/278857_percepttesting_sequence_01.csv"
$ws.Range("A4").Value = "sequences/278857_percepttesting_sequence_02.csv"
$ws.Range("A5").Value = "sequences/278857_percepttesting_sequence_03.csv"
$ws.Range("A6").Value = "sequences/278857_percepttesting_sequence_04.csv"
